# Add 2022-Q4 data:
#  - new "2022-Q4" worksheet inserted right after "总计" (pushing
#    "2022-Q1" / "2021-Q4" one slot to the right)
#  - a new row is added at the top of the "总计" summary table for the
#    2022-Q4 figures, with the two pre-existing rows re-written one row
#    lower (same values, new row numbers)

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) New "2022-Q4" sheet, positioned right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# header row
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

$hdr = $q4.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# data rows
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "'011527"
$q4.Cells.Item(2, 3).Value = "博时恒悦6个月持有期混合A"
$q4.Cells.Item(2, 4).Value = "'6.50"
$q4.Cells.Item(2, 5).Value = "'21.52"
$q4.Cells.Item(2, 6).Value = "'1.59"
$q4.Cells.Item(2, 7).Value = "'0.1034"
$q4.Cells.Item(2, 8).Value = 7

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "'002409"
$q4.Cells.Item(3, 3).Value = "华夏新活力灵活配置混合A"
$q4.Cells.Item(3, 4).Value = "'0.12"
$q4.Cells.Item(3, 5).Value = "'77.58"
$q4.Cells.Item(3, 6).Value = "'6.91"
$q4.Cells.Item(3, 7).Value = "'0.0083"
$q4.Cells.Item(3, 8).Value = 5

$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "'011528"
$q4.Cells.Item(4, 3).Value = "博时恒悦6个月持有期混合C"
$q4.Cells.Item(4, 4).Value = "'0.22"
$q4.Cells.Item(4, 5).Value = "'21.52"
$q4.Cells.Item(4, 6).Value = "'1.59"
$q4.Cells.Item(4, 7).Value = "'0.0035"
$q4.Cells.Item(4, 8).Value = 7

$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "'002410"
$q4.Cells.Item(5, 3).Value = "华夏新活力灵活配置混合C"
$q4.Cells.Item(5, 4).Value = "'0.00"
$q4.Cells.Item(5, 5).Value = "'77.58"
$q4.Cells.Item(5, 6).Value = "'6.91"
$q4.Cells.Item(5, 7).Value = 0
$q4.Cells.Item(5, 8).Value = 5

$col_a = $q4.Range("A2:A5")
$col_a.Font.Bold = $true
$col_a.HorizontalAlignment = -4108
$col_a.VerticalAlignment = -4160
$col_a.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2) "总计" sheet: add the 2022-Q4 row, push the other two rows down
# ---------------------------------------------------------------------
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q4"
$total.Cells.Item(4, 3).Value = 4
$total.Cells.Item(4, 4).Value = 0.42

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q1"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.45

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.12

$total.Range("A2").Copy()
$total.Range("A3:A4").PasteSpecial(-4122)
